# Update "想去人数" (want-to-go count) figures for two events that appear
# on both the "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 137
$wsExhibit.Range("F5").Value = 2912
$wsExhibit.Range("F7").Value = 395

# Sheet 4: 全部类型 (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 137
$wsAll.Range("F5").Value = 2912
$wsAll.Range("F9").Value = 395
